$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 "time_taken" -- copy the existing header style (bold,
# bordered, centered) from E1 so F1 matches the other header cells.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Per-row "time_taken" metadata values.
$ws.Range("F2").Value = "2021-10-05 10:50:13.003167"
$ws.Range("F3").Value = "2021-10-05 10:50:13.003183"
$ws.Range("F4").Value = "2021-10-05 10:50:13.003187"
$ws.Range("F5").Value = "2021-10-05 10:50:13.003190"
$ws.Range("F6").Value = "2021-10-05 10:50:13.003193"
